$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set A3 value to "x" (adds a new shared string)
$ws.Range("A3").Value = "x"

# Update the selected cell to C6
$ws.Range("C6").Select()
